# Automatische test-sync: 2025-06-17 21:55:47
# Appends a new incoming mail log row to the "Logs" sheet and bumps the
# matching "Informatieaanvraag" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 36

$logs.Cells.Item($newRow, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Cells.Item($newRow, 4).Value = "Informatieaanvraag"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nBedankt voor uw vraag. Onze openingstijden zijn maandag t/m vrijdag van 9:00 tot 17:00 uur. Op zaterdag zijn wij geopend van 10:00 tot 16:00 uur. Op zondag zijn wij gesloten.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-17 21:55:23"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# Keep the new row's height in line with the rest of the log (the multi-line
# reply text would otherwise trigger an autofit taller than its neighbours).
$logs.Rows.Item($newRow).RowHeight = $logs.Rows.Item($newRow - 1).RowHeight

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 17

# Extend the conditional-formatting ranges to cover the newly added row.
$categoryFormats = $logs.Range("D2:D35").FormatConditions
for ($i = 1; $i -le $categoryFormats.Count; $i++) {
    $categoryFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D36"))
}

$answeredFormats = $logs.Range("G2:G35").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G36"))
}
